$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.2245614035087719
$ws.Range("C2").Value = 0.4842105263157895
$ws.Range("J2").Value = 0.01754385964912281
$ws.Range("P2").Value = 0.1649122807017544
$ws.Range("S2").Value = 0.1087719298245614
$ws.Range("B3").Value = 0.006944444444444444
$ws.Range("C3").Value = 0.02083333333333333
$ws.Range("J3").Value = 0.04166666666666666
$ws.Range("P3").Value = 0.6736111111111112
$ws.Range("S3").Value = 0.2569444444444444
$ws.Range("P4").Value = 0.6578947368421053
$ws.Range("S4").Value = 0.3421052631578947
$ws.Range("B6").Value = 0.04477611940298507
$ws.Range("D6").Value = 0.01865671641791045
$ws.Range("F6").Value = 0.07462686567164178
$ws.Range("J6").Value = 0.2313432835820896
$ws.Range("O6").Value = 0.02985074626865672
$ws.Range("Q6").Value = 0.1417910447761194
$ws.Range("R6").Value = 0.04477611940298507
$ws.Range("S6").Value = 0.4141791044776119
$ws.Range("B7").Value = 0.108
$ws.Range("D7").Value = 0.012
$ws.Range("F7").Value = 0.06
$ws.Range("J7").Value = 0.16
$ws.Range("O7").Value = 0.044
$ws.Range("Q7").Value = 0.112
$ws.Range("R7").Value = 0.056
$ws.Range("S7").Value = 0.448
$ws.Range("B8").Value = 0.08151093439363817
$ws.Range("D8").Value = 0.007952286282306162
$ws.Range("E8").Value = 0.001988071570576541
$ws.Range("F8").Value = 0.08349900596421471
$ws.Range("J8").Value = 0.08548707753479125
$ws.Range("O8").Value = 0.02584493041749503
$ws.Range("Q8").Value = 0.1292246520874752
$ws.Range("R8").Value = 0.07753479125248509
$ws.Range("S8").Value = 0.5069582504970179
$ws.Range("B9").Value = 0.04854368932038835
$ws.Range("D9").Value = 0.009708737864077669
$ws.Range("F9").Value = 0.0825242718446602
$ws.Range("J9").Value = 0.116504854368932
$ws.Range("O9").Value = 0.02427184466019417
$ws.Range("Q9").Value = 0.1359223300970874
$ws.Range("R9").Value = 0.04368932038834952
$ws.Range("S9").Value = 0.5388349514563107
$ws.Range("B10").Value = 0.1128608923884514
$ws.Range("D10").Value = 0.02099737532808399
$ws.Range("F10").Value = 0.09011373578302712
$ws.Range("J10").Value = 0.1321084864391951
$ws.Range("O10").Value = 0.01924759405074366
$ws.Range("Q10").Value = 0.1881014873140857
$ws.Range("R10").Value = 0.03849518810148731
$ws.Range("S10").Value = 0.3980752405949257
$ws.Range("G11").Value = 0.1508951406649616
$ws.Range("J11").Value = 0.06649616368286446
$ws.Range("K11").Value = 0.1969309462915601
$ws.Range("L11").Value = 0.5473145780051151
$ws.Range("S11").Value = 0.03836317135549872
$ws.Range("G12").Value = 0.7359307359307359
$ws.Range("J12").Value = 0.1904761904761905
$ws.Range("K12").Value = 0.01298701298701299
$ws.Range("L12").Value = 0.02164502164502164
$ws.Range("S12").Value = 0.03896103896103896
$ws.Range("G13").Value = 0.6875
$ws.Range("J13").Value = 0.1875
$ws.Range("S13").Value = 0.125
$ws.Range("F15").Value = 0.03703703703703703
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.06481481481481481
$ws.Range("J15").Value = 0.2731481481481481
$ws.Range("K15").Value = 0.08333333333333333
$ws.Range("N15").Value = 0.004629629629629629
$ws.Range("O15").Value = 0.06481481481481481
$ws.Range("S15").Value = 0.3055555555555556
$ws.Range("F16").Value = 0.03048780487804878
$ws.Range("H16").Value = 0.2195121951219512
$ws.Range("I16").Value = 0.06097560975609756
$ws.Range("J16").Value = 0.3353658536585366
$ws.Range("K16").Value = 0.1097560975609756
$ws.Range("M16").Value = 0.01829268292682927
$ws.Range("O16").Value = 0.0426829268292683
$ws.Range("S16").Value = 0.1829268292682927
$ws.Range("F17").Value = 0.01055408970976253
$ws.Range("H17").Value = 0.2058047493403694
$ws.Range("I17").Value = 0.1213720316622691
$ws.Range("J17").Value = 0.3614775725593667
$ws.Range("K17").Value = 0.1055408970976253
$ws.Range("M17").Value = 0.01319261213720317
$ws.Range("N17").Value = 0.002638522427440633
$ws.Range("O17").Value = 0.05277044854881267
$ws.Range("S17").Value = 0.1266490765171504
$ws.Range("F18").Value = 0.0170940170940171
$ws.Range("H18").Value = 0.2136752136752137
$ws.Range("I18").Value = 0.1025641025641026
$ws.Range("J18").Value = 0.3162393162393162
$ws.Range("K18").Value = 0.1538461538461539
$ws.Range("M18").Value = 0.02564102564102564
$ws.Range("O18").Value = 0.03418803418803419
$ws.Range("S18").Value = 0.1367521367521368
$ws.Range("F19").Value = 0.0195945945945946
$ws.Range("H19").Value = 0.2195945945945946
$ws.Range("I19").Value = 0.08378378378378379
$ws.Range("J19").Value = 0.3182432432432433
$ws.Range("K19").Value = 0.1378378378378378
$ws.Range("M19").Value = 0.02567567567567568
$ws.Range("S19").Value = 0.1385135135135135
Write-Host "Applied 106 cell updates"
